# Edit script for horarios-141-2026-01-13.xlsx
# Commit: "Horarios actualizados Linea 141 - 268"
# Updates the scraped bus-schedule data across the 3 worksheets:
#   Sheet 1 "LP1912"      -> timestamp/total bump, a few row corrections, 12 new rows (179-190)
#   Sheet 2 "LP1912-215"  -> timestamp bump only
#   Sheet 3 "6203-6173"   -> timestamp/total bump, 1 new row (37)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "LP1912"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Cells.Item(2, 1).Value = "Última actualización: 10:56:30"
$ws1.Cells.Item(3, 1).Value = "Total filas: 185"

$sheet1Rows = @(
    @{R=39; A="06:43:40"; B="06:46"; C="225_C ROCA-H SUR"; D=3; E="LP1912"},
    @{R=40; A="05:18:56"; B="06:46"; C="215C_EL PATO"; D=88; E="LP1912"},
    @{R=62; A="05:49:40"; B="07:32"; C="84_COLONIA URQUIZA-ESC 49"; D=103; E="LP1912"},
    @{R=63; A="06:15:04"; B="07:32"; C="11_ETCHEVERRY"; D=77; E="LP1912"},
    @{R=109; A="07:47:32"; B="09:22"; C="17_ROMERO"; D=95; E="LP1912"},
    @{R=110; A="07:59:28"; B="09:22"; C="16_SANTA ANA"; D=83; E="LP1912"},
    @{R=119; A="08:57:13"; B="09:35"; C="16_SANTA ANA"; D=38; E="LP1912"},
    @{R=120; A="08:57:13"; B="09:35"; C="23_HERNANDEZ"; D=38; E="LP1912"},
    @{R=121; A="09:38:09"; B="09:41"; C="23_HERNANDEZ"; D=3; E="LP1912"},
    @{R=122; A="08:21:50"; B="09:41"; C="215C_EL PATO"; D=80; E="LP1912"},
    @{R=123; A="09:38:09"; B="09:41"; C="14_ABASTO"; D=3; E="LP1912"},
    @{R=155; A="10:56:30"; B="10:57"; C="16_SANTA ANA"; D=1; E="LP1912"},
    @{R=156; A="09:38:09"; B="10:58"; C="27_EL RETIRO"; D=80; E="LP1912"},
    @{R=157; A="10:26:41"; B="11:01"; C="215C_EL PATO"; D=35; E="LP1912"},
    @{R=158; A="09:38:09"; B="11:02"; C="215C_EL PATO"; D=84; E="LP1912"},
    @{R=159; A="10:26:41"; B="11:03"; C="11_ETCHEVERRY"; D=37; E="LP1912"},
    @{R=160; A="10:26:41"; B="11:04"; C="23_HERNANDEZ"; D=38; E="LP1912"},
    @{R=161; A="10:26:41"; B="11:06"; C="16_P MOR-167 Y 521"; D=40; E="LP1912"},
    @{R=162; A="09:38:09"; B="11:07"; C="16_P MOR-167 Y 521"; D=89; E="LP1912"},
    @{R=163; A="10:56:30"; B="11:08"; C="23_HERNANDEZ"; D=12; E="LP1912"},
    @{R=164; A="10:26:41"; B="11:12"; C="15_ABASTO"; D=46; E="LP1912"},
    @{R=165; A="10:26:41"; B="11:19"; C="86_EST CHICA-ESC AGRARIA"; D=53; E="LP1912"},
    @{R=166; A="09:38:09"; B="11:20"; C="86_EST CHICA-ESC AGRARIA"; D=102; E="LP1912"},
    @{R=167; A="09:38:09"; B="11:21"; C="26_HERNANDEZ"; D=103; E="LP1912"},
    @{R=168; A="10:56:30"; B="11:24"; C="10_OLMOS"; D=28; E="LP1912"},
    @{R=169; A="09:38:09"; B="11:27"; C="225_C ROCA-H SUR"; D=109; E="LP1912"},
    @{R=170; A="09:38:09"; B="11:32"; C="81_EL PELIGRO"; D=114; E="LP1912"},
    @{R=171; A="10:56:30"; B="11:34"; C="23_HERNANDEZ"; D=38; E="LP1912"},
    @{R=172; A="09:38:09"; B="11:35"; C="11_ETCHEVERRY"; D=69; E="LP1912"},
    @{R=173; A="09:38:09"; B="11:36"; C="11_ETCHEVERRY"; D=118; E="LP1912"},
    @{R=174; A="10:26:41"; B="11:41"; C="17_ROMERO"; D=75; E="LP1912"},
    @{R=175; A="10:56:30"; B="11:42"; C="215B_EL PATO"; D=46; E="LP1912"},
    @{R=176; A="09:38:09"; B="11:51"; C="215B_EL PATO"; D=85; E="LP1912"},
    @{R=177; A="10:56:30"; B="11:52"; C="15_ABASTO"; D=56; E="LP1912"},
    @{R=178; A="09:38:09"; B="11:59"; C="225_GOMEZ"; D=93; E="LP1912"},
    @{R=179; A="10:26:41"; B="12:02"; C="84_COLONIA URQUIZA-ESC 49"; D=96; E="LP1912"},
    @{R=180; A="10:26:41"; B="12:06"; C="16_P MOR-SANTA ANA"; D=100; E="LP1912"},
    @{R=181; A="10:56:30"; B="12:06"; C="14_ABASTO"; D=70; E="LP1912"},
    @{R=182; A="10:56:30"; B="12:10"; C="10_OLMOS"; D=74; E="LP1912"},
    @{R=183; A="10:26:41"; B="12:14"; C="17_ROMERO"; D=108; E="LP1912"},
    @{R=184; A="10:26:41"; B="12:19"; C="14_ABASTO"; D=113; E="LP1912"},
    @{R=185; A="10:26:41"; B="12:20"; C="215A_EL PATO"; D=114; E="LP1912"},
    @{R=186; A="10:56:30"; B="12:20"; C="14_ABASTO"; D=84; E="LP1912"},
    @{R=187; A="10:26:41"; B="12:21"; C="26_HERNANDEZ"; D=115; E="LP1912"},
    @{R=188; A="10:56:30"; B="12:36"; C="27_EL RETIRO"; D=100; E="LP1912"},
    @{R=189; A="10:56:30"; B="12:38"; C="17_179 Y 38"; D=102; E="LP1912"},
    @{R=190; A="10:56:30"; B="12:41"; C="10_OLMOS"; D=105; E="LP1912"}
)

foreach ($row in $sheet1Rows) {
    $ws1.Cells.Item($row.R, 1).Value = $row.A
    $ws1.Cells.Item($row.R, 2).Value = $row.B
    $ws1.Cells.Item($row.R, 3).Value = $row.C
    $ws1.Cells.Item($row.R, 4).Value = $row.D
    $ws1.Cells.Item($row.R, 5).Value = $row.E
}

# ---------------------------------------------------------------------------
# Sheet 2: "LP1912-215"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Cells.Item(2, 1).Value = "Última actualización: 10:56:30"

# ---------------------------------------------------------------------------
# Sheet 3: "6203-6173"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Cells.Item(2, 1).Value = "Última actualización: 10:56:30"
$ws3.Cells.Item(3, 1).Value = "Total filas: 32"

$ws3.Cells.Item(37, 1).Value = "10:56:30"
$ws3.Cells.Item(37, 2).Value = "12:54"
$ws3.Cells.Item(37, 3).Value = "215C_LA PLATA"
$ws3.Cells.Item(37, 4).Value = 118
$ws3.Cells.Item(37, 5).Value = "L6203"
